$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 125.583336
$ws.Range("I33").Value = 85.28570999999999
$ws.Range("J33").Value = 182
$ws.Range("K33").Value = 85.28570999999999
$ws.Range("L33").Value = 182
$ws.Range("M33").Value = 143.71429
$ws.Range("N33").Value = -640
# Row 40
$ws.Range("H40").Value = 5744.595
$ws.Range("I40").Value = 3553.6538
$ws.Range("K40").Value = 3553.6538
$ws.Range("M40").Value = -3378.6538
# Row 64
$ws.Range("H64").Value = 9999
$ws.Range("I64").Value = 9999
$ws.Range("K64").Value = 9999
$ws.Range("M64").Value = -9751
# Row 67
$ws.Range("H67").Value = 9999
$ws.Range("I67").Value = 9999
$ws.Range("K67").Value = 9999
$ws.Range("M67").Value = -9141
# Row 93
$ws.Range("H93").Value = 60601
$ws.Range("J93").Value = 60601
$ws.Range("L93").Value = 60601
$ws.Range("N93").Value = -65593

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 795.58826
$ws.Range("I2").Value = 594.7143
$ws.Range("J2").Value = 1733
$ws.Range("K2").Value = 594.7143
$ws.Range("L2").Value = 1733
$ws.Range("M2").Value = -481.7143
$ws.Range("N2").Value = -1959
# Row 32
$ws.Range("H32").Value = 2818.543
$ws.Range("I32").Value = 2686.3635
$ws.Range("K32").Value = 2686.3635
$ws.Range("M32").Value = -2399.3635
# Row 116
$ws.Range("H116").Value = 795.58826
$ws.Range("I116").Value = 594.7143
$ws.Range("J116").Value = 1733
$ws.Range("K116").Value = 594.7143
$ws.Range("L116").Value = 1733
$ws.Range("M116").Value = 1699.2857
$ws.Range("N116").Value = -6321

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 795.58826
$ws.Range("I3").Value = 594.7143
$ws.Range("J3").Value = 1733
$ws.Range("K3").Value = 594.7143
$ws.Range("L3").Value = 1733
$ws.Range("M3").Value = -480.7143
$ws.Range("N3").Value = -1961
# Row 8
$ws.Range("H8").Value = 801.3333
$ws.Range("I8").Value = 801.3333
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 801.3333
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -661.3333
$ws.Range("N8").Value = ""
# Row 86
$ws.Range("H86").Value = 4970.9287
$ws.Range("I86").Value = 3425.625
$ws.Range("J86").Value = 7031.3335
$ws.Range("K86").Value = 3425.625
$ws.Range("L86").Value = 7031.3335
$ws.Range("M86").Value = -2302.625
$ws.Range("N86").Value = -9277.333500000001
# Row 89
$ws.Range("H89").Value = 4970.9287
$ws.Range("I89").Value = 3425.625
$ws.Range("J89").Value = 7031.3335
$ws.Range("K89").Value = 17128.125
$ws.Range("L89").Value = 35156.6675
$ws.Range("M89").Value = -11512.125
$ws.Range("N89").Value = -46388.6675

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6093.625
$ws.Range("I31").Value = 2833.8
$ws.Range("J31").Value = 8422.071
$ws.Range("K31").Value = 2833.8
$ws.Range("L31").Value = 8422.071
$ws.Range("M31").Value = -2538.8
$ws.Range("N31").Value = -9012.071
# Row 34
$ws.Range("H34").Value = 6093.625
$ws.Range("I34").Value = 2833.8
$ws.Range("J34").Value = 8422.071
$ws.Range("K34").Value = 2833.8
$ws.Range("L34").Value = 8422.071
$ws.Range("M34").Value = -2631.8
$ws.Range("N34").Value = -8826.071
# Row 58
$ws.Range("H58").Value = 3151.1052
$ws.Range("I58").Value = 2680.6875
$ws.Range("J58").Value = 5660
$ws.Range("K58").Value = 2680.6875
$ws.Range("L58").Value = 5660
$ws.Range("M58").Value = -2477.6875
$ws.Range("N58").Value = -6066
# Row 105
$ws.Range("H105").Value = 2237.5
$ws.Range("I105").Value = 1569.5
$ws.Range("K105").Value = 1569.5
$ws.Range("M105").Value = 177.5
# Row 107
$ws.Range("H107").Value = 182.33333
$ws.Range("I107").Value = 118.8
$ws.Range("K107").Value = 118.8
$ws.Range("M107").Value = 1801.2
# Row 122
$ws.Range("H122").Value = 1391.8182
$ws.Range("I122").Value = 1391.8182
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4175.4546
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1725.4546
$ws.Range("N122").Value = ""
# Row 132
$ws.Range("H132").Value = 2827.5
$ws.Range("I132").Value = 2367.1428
$ws.Range("K132").Value = 7101.428400000001
$ws.Range("M132").Value = -4571.428400000001
# Row 134
$ws.Range("H134").Value = 1798.0714
$ws.Range("J134").Value = 1100
$ws.Range("L134").Value = 3300
$ws.Range("N134").Value = -8370
# Row 136
$ws.Range("H136").Value = 3151.1052
$ws.Range("I136").Value = 2680.6875
$ws.Range("J136").Value = 5660
$ws.Range("K136").Value = 8042.0625
$ws.Range("L136").Value = 16980
$ws.Range("M136").Value = -5492.0625
$ws.Range("N136").Value = -22080

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 949.1111
$ws.Range("J34").Value = 1399.8334
$ws.Range("L34").Value = 4199.5002
$ws.Range("N34").Value = -4367.5002
# Row 39
$ws.Range("H39").Value = 5425
$ws.Range("J39").Value = 7980
$ws.Range("L39").Value = 23940
$ws.Range("N39").Value = -24528
# Row 50
$ws.Range("H50").Value = 464.25
$ws.Range("I50").Value = 464.25
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 1392.75
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -911.75
$ws.Range("N50").Value = ""
# Row 53
$ws.Range("H53").Value = 464.25
$ws.Range("I53").Value = 464.25
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 1392.75
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -911.75
$ws.Range("N53").Value = ""
# Row 55
$ws.Range("H55").Value = 5121.8
$ws.Range("J55").Value = 6151.875
$ws.Range("L55").Value = 18455.625
$ws.Range("N55").Value = -18809.625
# Row 116
$ws.Range("H116").Value = 2617
$ws.Range("J116").Value = 3261
$ws.Range("L116").Value = 9783
$ws.Range("N116").Value = -16667
# Row 125
$ws.Range("H125").Value = 9990
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 9990
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 29970
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = -39810

$ws = $wb.Worksheets.Item("GSM")
# Row 95
$ws.Range("H95").Value = 23997.25
$ws.Range("J95").Value = 23997.25
$ws.Range("L95").Value = 23997.25
$ws.Range("N95").Value = -29489.25
# Row 97
$ws.Range("H97").Value = 636.7778
$ws.Range("I97").Value = 604
$ws.Range("J97").Value = 899
$ws.Range("K97").Value = 604
$ws.Range("L97").Value = 899
$ws.Range("M97").Value = -108
$ws.Range("N97").Value = -1891
# Row 122
$ws.Range("H122").Value = 2592.5
$ws.Range("I122").Value = 1768.625
$ws.Range("J122").Value = 5888
$ws.Range("K122").Value = 5305.875
$ws.Range("L122").Value = 17664
$ws.Range("M122").Value = -2855.875
$ws.Range("N122").Value = -22564

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 5708
$ws.Range("J46").Value = 6596.6
$ws.Range("L46").Value = 6596.6
$ws.Range("N46").Value = -6972.6
# Row 61
$ws.Range("H61").Value = 3774.25
$ws.Range("I61").Value = 965.3333
$ws.Range("K61").Value = 965.3333
$ws.Range("M61").Value = -763.3333
# Row 93
$ws.Range("H93").Value = 901
$ws.Range("I93").Value = 851.5
$ws.Range("K93").Value = 851.5
$ws.Range("M93").Value = 396.5
# Row 100
$ws.Range("H100").Value = 5169.857
$ws.Range("I100").Value = 2486.4443
$ws.Range("K100").Value = 2486.4443
$ws.Range("M100").Value = -1945.4443
# Row 113
$ws.Range("H113").Value = 3774.25
$ws.Range("I113").Value = 965.3333
$ws.Range("K113").Value = 965.3333
$ws.Range("M113").Value = 1204.6667

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2195.5715
$ws.Range("I122").Value = 1047.5
$ws.Range("J122").Value = 5869.4
$ws.Range("K122").Value = 3142.5
$ws.Range("L122").Value = 17608.2
$ws.Range("M122").Value = -692.5
$ws.Range("N122").Value = -22508.2
# Row 123
$ws.Range("H123").Value = 48000
$ws.Range("J123").Value = 48000
$ws.Range("L123").Value = 48000
$ws.Range("N123").Value = -57800
# Row 126
$ws.Range("H126").Value = 4953.75
$ws.Range("I126").Value = 3042.1428
$ws.Range("J126").Value = 7630
$ws.Range("K126").Value = 9126.428400000001
$ws.Range("L126").Value = 22890
$ws.Range("M126").Value = -6656.428400000001
$ws.Range("N126").Value = -27830
